# Update "想去人数" (interest/attendance counter) values in column F
# for the 展览 (Exhibition) and 全部类型 (All Types) sheets, reflecting the
# latest scrape used to regenerate the gh-pages output.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# 展览 sheet (row -> new value)
$wsExhibition.Range("F6").Value  = 559
$wsExhibition.Range("F7").Value  = 1750
$wsExhibition.Range("F10").Value = 140
$wsExhibition.Range("F11").Value = 1923
$wsExhibition.Range("F12").Value = 138
$wsExhibition.Range("F13").Value = 253
$wsExhibition.Range("F14").Value = 441
$wsExhibition.Range("F15").Value = 11
$wsExhibition.Range("F16").Value = 278
$wsExhibition.Range("F21").Value = 41
$wsExhibition.Range("F22").Value = 54
$wsExhibition.Range("F23").Value = 1036
$wsExhibition.Range("F25").Value = 320
$wsExhibition.Range("F27").Value = 257
$wsExhibition.Range("F28").Value = 288

# 全部类型 sheet (row -> new value)
$wsAllTypes.Range("F6").Value  = 559
$wsAllTypes.Range("F7").Value  = 1750
$wsAllTypes.Range("F11").Value = 140
$wsAllTypes.Range("F12").Value = 1923
$wsAllTypes.Range("F13").Value = 138
$wsAllTypes.Range("F14").Value = 253
$wsAllTypes.Range("F15").Value = 441
$wsAllTypes.Range("F16").Value = 11
$wsAllTypes.Range("F17").Value = 278
$wsAllTypes.Range("F22").Value = 41
$wsAllTypes.Range("F23").Value = 54
$wsAllTypes.Range("F24").Value = 1036
$wsAllTypes.Range("F26").Value = 320
$wsAllTypes.Range("F28").Value = 257
$wsAllTypes.Range("F29").Value = 288
